# Update comp_percentile_rank_dog output: refresh HealthyDistance (AK) and
# TotalRiskScore (AM) percentile values now that InsertDB has been removed
# from the calculation code.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AK3").Value = 18.1
$ws.Range("AM3").Value = 44.6
$ws.Range("AM8").Value = 83.8
$ws.Range("AK10").Value = 38.7
$ws.Range("AM10").Value = 15.7
$ws.Range("AK12").Value = 35.8
$ws.Range("AK15").Value = 16.2
$ws.Range("AM15").Value = 36.8
$ws.Range("AM16").Value = 52.5
$ws.Range("AM17").Value = 82.40000000000001
$ws.Range("AK25").Value = 46.6
$ws.Range("AK26").Value = 87.7
$ws.Range("AM26").Value = 91.7
$ws.Range("AM27").Value = 65.2
$ws.Range("AM30").Value = 60.3
$ws.Range("AK38").Value = 57.4
$ws.Range("AK40").Value = 49.5
$ws.Range("AM40").Value = 71.09999999999999
$ws.Range("AM42").Value = 85.8
$ws.Range("AK56").Value = 36.8
$ws.Range("AK59").Value = 71.59999999999999
$ws.Range("AM59").Value = 43.1
$ws.Range("AK60").Value = 75
$ws.Range("AM60").Value = 50.5
$ws.Range("AK61").Value = 41.7
$ws.Range("AM61").Value = 17.2
$ws.Range("AM62").Value = 86.8
$ws.Range("AM64").Value = 84.8
$ws.Range("AK65").Value = 75.5
$ws.Range("AK67").Value = 55.4
$ws.Range("AM67").Value = 33.3
$ws.Range("AK68").Value = 60.3
$ws.Range("AM68").Value = 23
$ws.Range("AM69").Value = 67.2
$ws.Range("AK70").Value = 23
$ws.Range("AM70").Value = 18.6
$ws.Range("AK71").Value = 12.3
$ws.Range("AM71").Value = 11.8
$ws.Range("AM73").Value = 73
$ws.Range("AK75").Value = 80.40000000000001
$ws.Range("AM75").Value = 69.09999999999999
$ws.Range("AK77").Value = 21.1
$ws.Range("AK86").Value = 52.5
$ws.Range("AM86").Value = 63.2
$ws.Range("AK90").Value = 33.8
$ws.Range("AK95").Value = 78.90000000000001
$ws.Range("AM95").Value = 54.4
